$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Percent Complete" (E) and "Has Description" (G) values for rows 19, 20, 22, 23
$ws.Range("E19").Value = 10
$ws.Range("G19").Value = 1

$ws.Range("E20").Value = 10

$ws.Range("E22").Value = 10
$ws.Range("G22").Value = 1

$ws.Range("E23").Value = 10
$ws.Range("G23").Value = 1

# Force recalculation of totals row formulas
$excel.Calculate()

# Update selected cell to match the author's final cursor position
$ws.Range("G24").Select()
